# Added WAV Trigger software
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-apply the existing date/number formats (style indices 1 and 2) to the
# new cells by copying formats from already-styled cells, so the engine
# reuses the existing style entries instead of creating new ones.
$ws.Range("C3:D3").Copy()
$ws.Range("C4:D7").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 4 - WAV Trigger Firmware
$ws.Range("A4").Value = "WAV Trigger Firmware"
$ws.Range("B4").Value = "WAV Trigger"
$ws.Range("C4").Value = 41746
$ws.Range("D4").Value = 0.61
$ws.Range("E4").Value = "http://robertsonics.com/wav-trigger-downloads/"
$ws.Range("F4").Value = "Firmware for WAV Trigger"

# Row 5 - WAV Trigger Firmware Update Utility
$ws.Range("A5").Value = "WAV Trigger Firmware Update Utility"
$ws.Range("B5").Value = "WAV Trigger"
$ws.Range("C5").Value = 41746
$ws.Range("D5").Value = 1.02
$ws.Range("E5").Value = "http://robertsonics.com/wav-trigger-downloads/"
$ws.Range("F5").Value = "Firmware updater for WAV Trigger"

# Row 6 - WAV Trigger Init File Maker Utility
$ws.Range("A6").Value = "WAV Trigger Init File Maker Utility"
$ws.Range("B6").Value = "WAV Trigger"
$ws.Range("C6").Value = 41746
$ws.Range("D6").Value = 1.02
$ws.Range("E6").Value = "http://robertsonics.com/wav-trigger-downloads/"
$ws.Range("F6").Value = "Init File Maker for WAV Trigger"

# Row 7 - WAV Trigger Serial Control Utility
$ws.Range("A7").Value = "WAV Trigger Serial Control Utility"
$ws.Range("B7").Value = "WAV Trigger"
$ws.Range("C7").Value = 41746
$ws.Range("D7").Value = 0.23
$ws.Range("E7").Value = "http://robertsonics.com/wav-trigger-downloads/"
$ws.Range("F7").Value = "Serial Control Utility for WAV Trigger"

# Row 8 - WAV Trigger Online Guide (no version/D value)
$ws.Range("A8").Value = "WAV Trigger Online Guide"
$ws.Range("B8").Value = "WAV Trigger"
$ws.Range("C8").Value = 41746
$ws.Range("E8").Value = "http://robertsonics.com/wav-trigger-online-user-guide/"
$ws.Range("F8").Value = "WAV Trigger Online Guide"

# Update column widths to better fit the new, wider content (values chosen
# so the resulting stored width is as close as possible to the authored
# widths of 34.42578125 / 16.28515625 / 7.85546875 / 52.28515625 / 34.140625).
$ws.Columns.Item(1).ColumnWidth = 33.59375
$ws.Columns.Item(3).ColumnWidth = 15.4296875
$ws.Columns.Item(4).ColumnWidth = 7.03125
$ws.Columns.Item(5).ColumnWidth = 51.5625
$ws.Columns.Item(6).ColumnWidth = 33.3984375

# Update the active selection to mirror the authored workbook state.
$ws.Range("E15").Select() | Out-Null

Write-Host "WAV Trigger rows added"
